$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = $ws.Range("A2").Value2
$ws.Range("B3").Value = $ws.Range("B2").Value2
$ws.Range("C3").Value = "eligitel@gmail.com"
$ws.Range("D3").Value = $ws.Range("C2").Value2
$ws.Range("E3").Value = $ws.Range("E2").Value2
$ws.Range("F3").Value = "great information about bitcoin. Really helpful!"

$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:eligitel@gmail.com", "", "", "eligitel@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:ronenchen27@gmail.com", "", "", "ronenchen27@gmail.com")

$ws.Range("A2:F2").Copy()
$ws.Range("A3:F3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("F3").Select() | Out-Null
